$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the duplicated "value" header cells in C1:D1:E1:F1 (only A1/B1 remain)
$ws.Range("C1:F1").ClearContents()

# Row 8: "Model" -> "production_function" (value in B8 stays "Sigmoid")
$ws.Range("A8").Value2 = "production_function"

# Insert a new row 9 for the L_curve parameter (value 0, scientific number format
# matching the other numeric optimization parameters above it)
$ws.Rows("9").Insert()
$ws.Range("A9").Value2 = "L_curve"
$ws.Range("B9").Value2 = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now row 17, after the insert above) is no longer used
# and gets dropped entirely
$ws.Rows("17").Delete()

# Make optimization_parameters the active/selected sheet, with C1:F1 selected
$ws.Activate()
$ws.Range("C1:F1").Select()
